$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values could be misinterpreted as numbers by Excel,
# so they are stored as text exactly like the source data (matches original inlineStr cells).
$textRefs = @("D4", "D5", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '25.021.94'
$ws.Range("E2").Value = '  -3.75%  '
$ws.Range("D3").Value = '1.648.53'
$ws.Range("E3").Value = '  -5.49%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '236.64'
$ws.Range("E5").Value = '  -5.65%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.4820'
$ws.Range("E7").Value = '  -6.51%  '
$ws.Range("D8").Value = '0.2615'
$ws.Range("E8").Value = '  -5.31%  '
$ws.Range("D9").Value = '0.06003'
$ws.Range("E9").Value = '  -3.17%  '
$ws.Range("D10").Value = '0.07193'
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("D11").Value = '1.646.02'
$ws.Range("E11").Value = '  -5.67%  '
$ws.Range("D12").Value = '14.80'
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("D13").Value = '0.6201'
$ws.Range("E13").Value = '  -4.99%  '
$ws.Range("D14").Value = '4.579'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '72.82'
$ws.Range("E15").Value = '  -6.53%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = '25.003.25'
$ws.Range("E18").Value = '  -3.90%  '
$ws.Range("D19").Value = '11.51'
$ws.Range("E19").Value = '  -3.01%  '
$ws.Range("D20").Value = '0.000006609'
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("D21").Value = '4.545'
$ws.Range("E21").Value = '  +5.51%  '
$ws.Range("D22").Value = '1.855.98'
$ws.Range("E22").Value = '  -5.68%  '
$ws.Range("D23").Value = '8.608'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").Value = '5.297'
$ws.Range("E24").Value = '  -1.85%  '
$ws.Range("D25").Value = '132.10'
$ws.Range("E25").Value = '  -3.00%  '
$ws.Range("D26").Value = '14.92'
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("D27").Value = '1.394'
$ws.Range("E27").Value = '  -7.75%  '
$ws.Range("D28").Value = '103.03'
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").Value = '1.668'
$ws.Range("E29").Value = '  -7.01%  '
$ws.Range("D30").Value = '3.755'
$ws.Range("E30").Value = '  -5.18%  '
$ws.Range("D31").Value = '0.07889'
$ws.Range("E31").Value = '  -4.40%  '
$ws.Range("D32").Value = '3.589'
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("D33").Value = '0.04577'
$ws.Range("E33").Value = '  -2.27%  '
$ws.Range("D34").Value = '2.595'
$ws.Range("E34").Value = '  -2.32%  '
$ws.Range("D35").Value = '0.9365'
$ws.Range("E35").Value = '  -6.56%  '
$ws.Range("D36").Value = '0.5770'
$ws.Range("E36").Value = '  -7.93%  '
$ws.Range("D37").Value = '2.601'
$ws.Range("E37").Value = '  -4.71%  '
$ws.Range("D38").Value = '0.01561'
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("E39").Value = '  +9.97%  '
$ws.Range("D40").Value = '1.0000'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -5.46%  '
$ws.Range("D42").Value = '98.14'
$ws.Range("E42").Value = '  -2.53%  '
$ws.Range("D43").Value = '0.3724'
$ws.Range("E43").Value = '  -4.29%  '
$ws.Range("D44").Value = '4.787'
$ws.Range("E44").Value = '  -4.74%  '
$ws.Range("D45").Value = '0.1141'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").Value = '6.123'
$ws.Range("E46").Value = '  -3.86%  '
$ws.Range("D47").Value = '0.05195'
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("D48").Value = '29.77'
$ws.Range("E48").Value = '  -3.41%  '
$ws.Range("B49").Value = 'TrueUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D49").Value = '1.002'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '50.53'
$ws.Range("E50").Value = '  -9.18%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").Value = '  -0.38%  '

# Restore default (Normal) style on the cells we temporarily forced to text format,
# so no stray style/number-format is left behind on the cell.
foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = "Normal"
}
